$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17, column A: was a Leetcode problem number (189); change it to the
# "GFG" text used by GeeksforGeeks-sourced problems (same as column A on
# rows 3, 6-9, 13, 15).
$ws.Range("A17").Value = "GFG"

# Add new row 18 for "Largest Element in Array" (GFG / Java / 15-Mar-2023),
# matching the layout of the existing GFG rows above it.
$ws.Range("A17").Copy($ws.Range("A18"))
$ws.Range("A18").Value = "GFG"

$ws.Range("B18").Value = "Largest Element in Array"
$ws.Range("B18").Style = "Normal"

$ws.Range("C18").Value = "Java"

$ws.Range("D18").Value = 45000
$ws.Range("D18").NumberFormat = "d-mmm-yy"

# Update the active selection / scroll position, as recorded in the sheet view
$null = $ws.Range("I14").Select()
$excel.ActiveWindow.ScrollRow = 4

Write-Host "done"
